# Suma_Tiempos.xlsx edit:
#  - Insert 4 new rows (Etapa 4..7) between the existing "Etapa 3" row (row 5)
#    and the "Total" row (old row 6, which becomes row 10).
#  - Fill in the new rows' "Columna E" (5th day) times.
#  - Extend the Total row's E-column SUM formula to cover the new rows.
#  - The second mini-table below (old rows 9-14) is pushed down to rows 13-18
#    automatically by the row insertion, formulas included.
#  - Update the sheet dimension / active selection to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert four blank rows above the old "Total" row (row 6). Each insert
# pushes the Total row (and everything below it) down by one, and the new
# blank row inherits the number formatting of the row above it (style 3 for
# column A, style 4 for B:L), same as Excel normally behaves on a manual
# row insert.
$ws.Rows.Item(6).Insert()
$ws.Rows.Item(6).Insert()
$ws.Rows.Item(6).Insert()
$ws.Rows.Item(6).Insert()

# Labels for the new stage rows.
$ws.Range("A6").Value = "Etapa 4"
$ws.Range("A7").Value = "Etapa 5"
$ws.Range("A8").Value = "Etapa 6"
$ws.Range("A9").Value = "Etapa 7"

# Times recorded for column E (the only populated data column for these
# new stages).
$ws.Range("E6").Value = 0.0766782407407407
$ws.Range("E7").Value = 0.0718634259259259
$ws.Range("E8").Value = 0.0725578703703704
$ws.Range("E9").Value = 0.0760300925925926

# The Total row (now row 10) needs its column-E sum expanded to include
# the newly inserted rows 6-9 (the other columns keep summing just 3:5).
$ws.Range("E10").Formula = "=SUM(E3:E9)"

# Match the workbook's last saved selection.
$ws.Range("E13").Select()
